$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'23.184.40"
$ws.Range("E2").Value = "'  -3.22%  "
$ws.Range("D3").Value = "'1.609.00"
$ws.Range("E3").Value = "'  -2.67%  "
$ws.Range("D4").Value = "'1.0000"
$ws.Range("E4").Value = "'  -0.11%  "
$ws.Range("D5").Value = "'0.9999"
$ws.Range("E5").Value = "'  -0.10%  "
$ws.Range("D6").Value = "'303.09"
$ws.Range("E6").Value = "'  -2.13%  "
$ws.Range("E7").Value = "'  -3.23%  "
$ws.Range("D8").Value = "'0.3672"
$ws.Range("E8").Value = "'  -4.34%  "
$ws.Range("D9").Value = "'49.19"
$ws.Range("E9").Value = "'  -4.08%  "
$ws.Range("D10").Value = "'1.000"
$ws.Range("E10").Value = "'  -0.08%  "
$ws.Range("D11").Value = "'1.276"
$ws.Range("E11").Value = "'  -5.85%  "
$ws.Range("D12").Value = "'0.08099"
$ws.Range("E13").Value = "'  -3.48%  "
$ws.Range("D14").Value = "'6.625"
$ws.Range("E14").Value = "'  -6.90%  "
$ws.Range("D15").Value = "'7.611"
$ws.Range("E15").Value = "'  -3.49%  "
$ws.Range("D16").Value = "'0.00001271"
$ws.Range("E16").Value = "'  -3.50%  "
$ws.Range("D17").Value = "'1.609.33"
$ws.Range("E17").Value = "'  -2.64%  "
$ws.Range("D18").Value = "'91.66"
$ws.Range("E18").Value = "'  -3.15%  "
$ws.Range("D19").Value = "'0.06791"
$ws.Range("E19").Value = "'  -2.91%  "
$ws.Range("D20").Value = "'18.40"
$ws.Range("E20").Value = "'  -6.87%  "
$ws.Range("D21").Value = "'6.598"
$ws.Range("E21").Value = "'  -4.76%  "
$ws.Range("D22").Value = "'0.9995"
$ws.Range("E22").Value = "'  -0.07%  "
$ws.Range("E23").Value = "'  -4.13%  "
$ws.Range("D24").Value = "'23.193.66"
$ws.Range("E24").Value = "'  -3.16%  "
$ws.Range("D25").Value = "'2.373"
$ws.Range("E25").Value = "'  -4.47%  "
$ws.Range("D26").Value = "'2.918"
$ws.Range("E26").Value = "'  -3.13%  "
$ws.Range("E27").Value = "'  -4.52%  "
$ws.Range("D28").Value = "'150.81"
$ws.Range("E28").Value = "'  -0.21%  "
$ws.Range("D29").Value = "'5.244"
$ws.Range("E29").Value = "'  -3.84%  "
$ws.Range("D30").Value = "'132.63"
$ws.Range("E30").Value = "'  -4.86%  "
$ws.Range("D31").Value = "'2.416"
$ws.Range("E31").Value = "'  -2.99%  "
$ws.Range("D32").Value = "'7.004"
$ws.Range("D33").Value = "'1.786.16"
$ws.Range("E33").Value = "'  -2.71%  "
$ws.Range("D34").Value = "'0.9829"
$ws.Range("E34").Value = "'  -5.91%  "
$ws.Range("D35").Value = "'0.07759"
$ws.Range("E35").Value = "'  -4.38%  "
$ws.Range("D36").Value = "'0.02794"
$ws.Range("E36").Value = "'  -5.93%  "
$ws.Range("D37").Value = "'6.309"
$ws.Range("E37").Value = "'  -6.74%  "
$ws.Range("D38").Value = "'0.2559"
$ws.Range("E38").Value = "'  -4.66%  "
$ws.Range("D39").Value = "'10.13"
$ws.Range("E39").Value = "'  -6.80%  "
$ws.Range("D40").Value = "'0.08874"
$ws.Range("E40").Value = "'  -2.94%  "
$ws.Range("D41").Value = "'1.399"
$ws.Range("E41").Value = "'  -2.00%  "
$ws.Range("D42").Value = "'0.7179"
$ws.Range("D43").Value = "'12.80"
$ws.Range("E43").Value = "'  -4.80%  "
$ws.Range("D44").Value = "'15.97"
$ws.Range("E44").Value = "'  -2.66%  "
$ws.Range("D45").Value = "'0.6635"
$ws.Range("E45").Value = "'  -4.49%  "
$ws.Range("D46").Value = "'2.306"
$ws.Range("E46").Value = "'  -6.12%  "
$ws.Range("D47").Value = "'0.9989"
$ws.Range("E48").Value = "'  -2.65%  "
$ws.Range("D49").Value = "'0.08021"
$ws.Range("E49").Value = "'  -3.08%  "
$ws.Range("D50").Value = "'131.88"
$ws.Range("E50").Value = "'  -2.06%  "
$ws.Range("D51").Value = "'1.172"
$ws.Range("E51").Value = "'  -2.79%  "
